# "unify the conception of DataNode, DataTable, Entity."
#
# Renames the two sheets to match the new naming scheme, re-heights a
# couple of header rows on the first sheet, and leaves the workbook with
# the second sheet ("DataTable") active/selected with cell H32 highlighted
# -- mirroring the state the workbook was saved in after the edit.

$wb = $excel.ActiveWorkbook

# --- Sheet renames -------------------------------------------------------
$wsDataNode  = $wb.Worksheets.Item(1)   # was "Property1"
$wsDataTable = $wb.Worksheets.Item(2)   # was "Record_Station"

$wsDataNode.Name  = "DataNode"
$wsDataTable.Name = "DataTable"

# --- Row-height tweaks on the DataNode sheet -----------------------------
$wsDataNode.Rows.Item(1).RowHeight = 27
$wsDataNode.Rows.Item(8).RowHeight = 54

# --- Final selection / active sheet --------------------------------------
# The workbook was left on the DataTable sheet with H32 selected.
$wsDataTable.Activate()
$wsDataTable.Range("H32").Select()
